# Update the build/version string across the workbook.
#
# Old version string:
#   mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)
# New version string:
#   Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)

$wb = $excel.ActiveWorkbook

$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- "About" sheet ---
$about = $wb.Worksheets.Item("About")

# A2: "Version: <version string>"
$about.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended citation text embedding the version string
$about.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Alardinskaya Coal Mine, Russia, M0767, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S ("build_version") holds the raw version string for data rows 2-18
for ($row = 2; $row -le 18; $row++) {
    $data.Cells.Item($row, 19).Value = $newVersion  # column S = 19
}
